$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Cells.Item(9, 2).Value = 2024
$ws.Cells.Item(9, 3).Value = 2025
$ws.Cells.Item(9, 4).Value = 2026
$ws.Cells.Item(9, 5).Value = 2027
$ws.Cells.Item(9, 6).Value = 2028
$ws.Cells.Item(9, 7).Value = 2029
$ws.Cells.Item(9, 8).Value = 2030
$ws.Cells.Item(9, 9).Value = 2031
$ws.Cells.Item(9, 10).Value = 2032
$ws.Cells.Item(9, 11).Value = 2033
$ws.Cells.Item(9, 12).Value = 2034

# Row 11
$ws.Cells.Item(11, 2).Value = 1.14
$ws.Cells.Item(11, 3).Value = 1.19
$ws.Cells.Item(11, 4).Value = 1.32
$ws.Cells.Item(11, 5).Value = 1.39
$ws.Cells.Item(11, 6).Value = 1.39
$ws.Cells.Item(11, 7).Value = 1.36
$ws.Cells.Item(11, 8).Value = 1.37
$ws.Cells.Item(11, 9).Value = 1.37
$ws.Cells.Item(11, 10).Value = 1.38
$ws.Cells.Item(11, 11).Value = 1.38
$ws.Cells.Item(11, 12).Value = 1.37

# Row 14
$ws.Cells.Item(14, 2).Value = 1.41
$ws.Cells.Item(14, 3).Value = 1.45
$ws.Cells.Item(14, 4).Value = 1.63
$ws.Cells.Item(14, 5).Value = 1.71
$ws.Cells.Item(14, 6).Value = 1.71
$ws.Cells.Item(14, 7).Value = 1.68
$ws.Cells.Item(14, 8).Value = 1.69
$ws.Cells.Item(14, 9).Value = 1.69
$ws.Cells.Item(14, 10).Value = 1.7
$ws.Cells.Item(14, 11).Value = 1.7
$ws.Cells.Item(14, 12).Value = 1.7

# Row 15
$ws.Cells.Item(15, 4).Value = 0.65
$ws.Cells.Item(15, 5).Value = 0.73
$ws.Cells.Item(15, 6).Value = 0.73
$ws.Cells.Item(15, 7).Value = 0.6899999999999999
$ws.Cells.Item(15, 8).Value = 0.7
$ws.Cells.Item(15, 9).Value = 0.7
$ws.Cells.Item(15, 10).Value = 0.7
$ws.Cells.Item(15, 11).Value = 0.7
$ws.Cells.Item(15, 12).Value = 0.7

# Row 17
$ws.Cells.Item(17, 2).Value = 0.82
$ws.Cells.Item(17, 3).Value = 0.82
$ws.Cells.Item(17, 4).Value = 0.99
$ws.Cells.Item(17, 5).Value = 0.98
$ws.Cells.Item(17, 6).Value = 0.98
$ws.Cells.Item(17, 7).Value = 0.99
$ws.Cells.Item(17, 8).Value = 0.99
$ws.Cells.Item(17, 9).Value = 0.99
$ws.Cells.Item(17, 10).Value = 1
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 1

# Row 20
$ws.Cells.Item(20, 2).Value = 1.27
$ws.Cells.Item(20, 3).Value = 1.31
$ws.Cells.Item(20, 4).Value = 1.37
$ws.Cells.Item(20, 5).Value = 1.44
$ws.Cells.Item(20, 6).Value = 1.44
$ws.Cells.Item(20, 7).Value = 1.42
$ws.Cells.Item(20, 8).Value = 1.43
$ws.Cells.Item(20, 9).Value = 1.43
$ws.Cells.Item(20, 10).Value = 1.44
$ws.Cells.Item(20, 11).Value = 1.44
$ws.Cells.Item(20, 12).Value = 1.43

# Row 21
$ws.Cells.Item(21, 2).Value = 0.89
$ws.Cells.Item(21, 3).Value = 0.95
$ws.Cells.Item(21, 4).Value = 1.21
$ws.Cells.Item(21, 5).Value = 1.31
$ws.Cells.Item(21, 6).Value = 1.31
$ws.Cells.Item(21, 7).Value = 1.25
$ws.Cells.Item(21, 8).Value = 1.26
$ws.Cells.Item(21, 9).Value = 1.26
$ws.Cells.Item(21, 10).Value = 1.27
$ws.Cells.Item(21, 11).Value = 1.27
$ws.Cells.Item(21, 12).Value = 1.26

# Row 23
$ws.Cells.Item(23, 2).Value = 0.38
$ws.Cells.Item(23, 3).Value = 0.37
$ws.Cells.Item(23, 4).Value = 0.16
$ws.Cells.Item(23, 5).Value = 0.13
$ws.Cells.Item(23, 6).Value = 0.13
$ws.Cells.Item(23, 7).Value = 0.17
$ws.Cells.Item(23, 8).Value = 0.17
$ws.Cells.Item(23, 10).Value = 0.17
$ws.Cells.Item(23, 12).Value = 0.18

# Row 25
$ws.Cells.Item(25, 2).Value = 0.86
$ws.Cells.Item(25, 3).Value = 0.8
$ws.Cells.Item(25, 4).Value = 0.8
$ws.Cells.Item(25, 5).Value = 0.75
$ws.Cells.Item(25, 6).Value = 0.75
$ws.Cells.Item(25, 7).Value = 0.74
$ws.Cells.Item(25, 8).Value = 0.73
$ws.Cells.Item(25, 9).Value = 0.73
$ws.Cells.Item(25, 10).Value = 0.73
$ws.Cells.Item(25, 11).Value = 0.73
$ws.Cells.Item(25, 12).Value = 0.73

# Row 26
$ws.Cells.Item(26, 2).Value = 0.11
$ws.Cells.Item(26, 3).Value = 0.1
$ws.Cells.Item(26, 4).Value = 0.1
$ws.Cells.Item(26, 6).Value = 0.05
$ws.Cells.Item(26, 7).Value = 0.04
$ws.Cells.Item(26, 8).Value = 0.03
$ws.Cells.Item(26, 9).Value = 0.03
$ws.Cells.Item(26, 10).Value = 0.03
$ws.Cells.Item(26, 11).Value = 0.03
$ws.Cells.Item(26, 12).Value = 0.03

# Row 28
$ws.Cells.Item(28, 2).Value = 0.01
$ws.Cells.Item(28, 3).Value = -0.05
$ws.Cells.Item(28, 4).Value = -0.46
$ws.Cells.Item(28, 5).Value = -0.47
$ws.Cells.Item(28, 6).Value = -0.47
$ws.Cells.Item(28, 7).Value = -0.45
$ws.Cells.Item(28, 8).Value = -0.45
$ws.Cells.Item(28, 9).Value = -0.45
$ws.Cells.Item(28, 10).Value = -0.45
$ws.Cells.Item(28, 11).Value = -0.45
$ws.Cells.Item(28, 12).Value = -0.45

# Row 29
$ws.Cells.Item(29, 2).Value = -1.13
$ws.Cells.Item(29, 3).Value = -1.24
$ws.Cells.Item(29, 4).Value = -1.78
$ws.Cells.Item(29, 5).Value = -1.86
$ws.Cells.Item(29, 6).Value = -1.86
$ws.Cells.Item(29, 7).Value = -1.81
$ws.Cells.Item(29, 8).Value = -1.82
$ws.Cells.Item(29, 9).Value = -1.83
$ws.Cells.Item(29, 10).Value = -1.83
$ws.Cells.Item(29, 11).Value = -1.83
$ws.Cells.Item(29, 12).Value = -1.83

# Row 32
$ws.Cells.Item(32, 2).Value = -0.01
$ws.Cells.Item(32, 4).Value = -0.2
$ws.Cells.Item(32, 5).Value = -0.2
$ws.Cells.Item(32, 6).Value = -0.2
$ws.Cells.Item(32, 7).Value = -0.19
$ws.Cells.Item(32, 8).Value = -0.19
$ws.Cells.Item(32, 9).Value = -0.19
$ws.Cells.Item(32, 10).Value = -0.19
$ws.Cells.Item(32, 11).Value = -0.19
$ws.Cells.Item(32, 12).Value = -0.19

# Row 33
$ws.Cells.Item(33, 2).Value = 0.05
$ws.Cells.Item(33, 3).Value = -0.1
$ws.Cells.Item(33, 4).Value = -0.92
$ws.Cells.Item(33, 6).Value = -0.93
$ws.Cells.Item(33, 7).Value = -0.9
$ws.Cells.Item(33, 8).Value = -0.9
$ws.Cells.Item(33, 9).Value = -0.9
$ws.Cells.Item(33, 11).Value = -0.91
$ws.Cells.Item(33, 12).Value = -0.91

# Row 35
$ws.Cells.Item(35, 2).Value = -0.07000000000000001
$ws.Cells.Item(35, 3).Value = 0.08
$ws.Cells.Item(35, 4).Value = 0.72
$ws.Cells.Item(35, 5).Value = 0.73
$ws.Cells.Item(35, 6).Value = 0.73
$ws.Cells.Item(35, 7).Value = 0.7
$ws.Cells.Item(35, 8).Value = 0.7
$ws.Cells.Item(35, 9).Value = 0.71
$ws.Cells.Item(35, 10).Value = 0.71
$ws.Cells.Item(35, 11).Value = 0.72
$ws.Cells.Item(35, 12).Value = 0.71

# Row 38
$ws.Cells.Item(38, 2).Value = 0.01
$ws.Cells.Item(38, 3).Value = -0.05
$ws.Cells.Item(38, 4).Value = -0.46
$ws.Cells.Item(38, 5).Value = -0.47
$ws.Cells.Item(38, 6).Value = -0.47
$ws.Cells.Item(38, 7).Value = -0.45
$ws.Cells.Item(38, 8).Value = -0.45
$ws.Cells.Item(38, 9).Value = -0.45
$ws.Cells.Item(38, 10).Value = -0.45
$ws.Cells.Item(38, 11).Value = -0.45
$ws.Cells.Item(38, 12).Value = -0.45

# Row 39
$ws.Cells.Item(39, 2).Value = 1.49
$ws.Cells.Item(39, 3).Value = 1.49
$ws.Cells.Item(39, 4).Value = 1.57
$ws.Cells.Item(39, 5).Value = 1.58
$ws.Cells.Item(39, 6).Value = 1.58
$ws.Cells.Item(39, 7).Value = 1.54
$ws.Cells.Item(39, 8).Value = 1.54
$ws.Cells.Item(39, 9).Value = 1.55
$ws.Cells.Item(39, 10).Value = 1.55
$ws.Cells.Item(39, 11).Value = 1.55
$ws.Cells.Item(39, 12).Value = 1.54

# Row 41
$ws.Cells.Item(41, 2).Value = -1.48
$ws.Cells.Item(41, 3).Value = -1.53
$ws.Cells.Item(41, 4).Value = -2.03
$ws.Cells.Item(41, 5).Value = -2.05
$ws.Cells.Item(41, 6).Value = -2.05
$ws.Cells.Item(41, 7).Value = -1.99
$ws.Cells.Item(41, 8).Value = -1.99
$ws.Cells.Item(41, 9).Value = -2
$ws.Cells.Item(41, 10).Value = -2
$ws.Cells.Item(41, 11).Value = -2
$ws.Cells.Item(41, 12).Value = -2
